$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -4
$ws.Range("F8").Value = -9
$ws.Range("F9").Value = 2
$ws.Range("F14").Value = -5
$ws.Range("F17").Value = 0
$ws.Range("F22").Value = -1
$ws.Range("F28").Value = 1
$ws.Range("F30").Value = -4
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 2
$ws.Range("F38").Value = -5
$ws.Range("F42").Value = -5
$ws.Range("F44").Value = -5
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = -1
$ws.Range("F48").Value = -3
$ws.Range("F49").Value = -4
$ws.Range("F53").Value = 0
$ws.Range("F54").Value = 2
$ws.Range("F56").Value = 5
$ws.Range("F58").Value = -6
$ws.Range("F62").Value = -9
$ws.Range("F63").Value = -3
$ws.Range("F64").Value = -5
$ws.Range("F65").Value = -3
